{"js": "// The author merged a duplicated \"conditions\" explanation that originally\n// appeared (slightly differently worded) in two consecutive paragraphs:\n//   Para A (\"Wenn der Spieler das Spiel startet, ...\") had:\n//     \"Die Bedingungen sind hier zum einen, dass das Wort mindestens zwei\n//      Zeichen hat und zum anderen, dass es keine Leerzeichen oder\n//      Sonderzeichen enth\u00e4lt. \"\n//   Para B (\"Zun\u00e4chst haben wir eine neue Klasse Spiel angelegt, ...\") had:\n//     \"Die zweite Bedingung ist, dass man keine Sonderzeichen und Zahlen\n//      verwenden darf.\"\n// The commit (\"Dopplung zusammengef\u00fcgt\" = \"duplication merged together\")\n// removes the sentence from Para A and folds its wording into Para B,\n// replacing Para B's old (duplicate/inaccurate) sentence.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldSentenceA =\n  \"Die Bedingungen sind hier zum einen, dass das Wort mindestens zwei Zeichen hat und zum anderen, dass es keine Leerzeichen oder Sonderzeichen enth\u00e4lt. \";\nconst oldSentenceB =\n  \"Die zweite Bedingung ist, dass man keine Sonderzeichen und Zahlen verwenden darf.\";\nconst newSentenceB =\n  \"Die zwei anderen Bedingungen sind, dass das Wort zum einen mindestens zwei Zeichen hat und zum anderen, dass es keine Leerzeichen oder Sonderzeichen enth\u00e4lt.\";\n\nlet paraA = null;\nlet paraB = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (paraA === null && t.indexOf(oldSentenceA) !== -1) {\n    paraA = paragraphs.items[i];\n  }\n  if (paraB === null && t.indexOf(oldSentenceB) !== -1) {\n    paraB = paragraphs.items[i];\n  }\n}\n\nif (paraA) {\n  const newTextA = paraA.text.replace(oldSentenceA, \"\");\n  paraA.getRange().insertText(newTextA, Word.InsertLocation.replace);\n}\n\nif (paraB) {\n  const newTextB = paraB.text.replace(oldSentenceB, newSentenceB);\n  paraB.getRange().insertText(newTextB, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The author merged a duplicated \"conditions\" explanation that originally\n# appeared (slightly differently worded) in two consecutive paragraphs:\n#   Para A (\"Wenn der Spieler das Spiel startet, ...\") had:\n#     \"Die Bedingungen sind hier zum einen, dass das Wort mindestens zwei\n#      Zeichen hat und zum anderen, dass es keine Leerzeichen oder\n#      Sonderzeichen enth\u00e4lt. \"\n#   Para B (\"Zun\u00e4chst haben wir eine neue Klasse Spiel angelegt, ...\") had:\n#     \"Die zweite Bedingung ist, dass man keine Sonderzeichen und Zahlen\n#      verwenden darf.\"\n# The commit (\"Dopplung zusammengef\u00fcgt\" = \"duplication merged together\")\n# removes the sentence from Para A and folds its wording into Para B,\n# replacing Para B's old (duplicate/inaccurate) sentence.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the duplicated sentence from the first paragraph.\n$find1 = $d.Content.Find\n$find1.Execute(\n    \"Die Bedingungen sind hier zum einen, dass das Wort mindestens zwei Zeichen hat und zum anderen, dass es keine Leerzeichen oder Sonderzeichen enth\u00e4lt. \",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"\",\n    2\n) | Out-Null\n\n# 2) Replace the sentence in the second paragraph with the reworded version\n#    that now carries the (merged) condition text.\n$find2 = $d.Content.Find\n$find2.Execute(\n    \"Die zweite Bedingung ist, dass man keine Sonderzeichen und Zahlen verwenden darf.\",\n    $true, $false, $false, $false, $false, $true, 1, $false,\n    \"Die zwei anderen Bedingungen sind, dass das Wort zum einen mindestens zwei Zeichen hat und zum anderen, dass es keine Leerzeichen oder Sonderzeichen enth\u00e4lt.\",\n    2\n) | Out-Null\n"}
